# "Roll back changes on Visium"
#
# A previous commit had added two lookup sheets/columns to the Visium
# template -- "preparation_instrument_vendor" (col Q) and
# "preparation_instrument_model" (col R) -- which pushed the existing
# ".metadata" bookkeeping sheet/column from Q to S. This change reverts
# that: the two preparation_instrument_* lookup sheets (and the Visium
# columns that validate against them) are removed, so ".metadata" slides
# back down to column Q. The .metadata sheet's recorded creation
# timestamp is also refreshed.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets("Visium")

# Column S currently carries the comment that documents the
# metadata_schema_id column; once Q and R are deleted, S becomes the new
# Q, so move that comment over (and drop the comments that belonged to
# the two columns being removed) before touching the layout.
$metadataComment = $ws.Range("S1").Comment.Text()
$null = $ws.Range("Q1").Comment.Text($metadataComment)
$ws.Range("R1").Comment.Delete()
$ws.Range("S1").Comment.Delete()

# Drop the preparation_instrument_vendor (Q) and preparation_instrument_model
# (R) columns from the main sheet -- this also removes their data
# validations and shifts the old S column (.metadata id) into Q.
$ws.Columns("Q:R").Delete()

# The two lookup sheets backing those columns are no longer referenced.
$wb.Worksheets("preparation_instrument_vendor").Delete()
$wb.Worksheets("preparation_instrument_model").Delete()

# Refresh the recorded creation timestamp on the .metadata sheet.
$meta = $wb.Worksheets(".metadata")
$meta.Cells.Item(2, 3).Value = "2023-11-01T15:37:30-07:00"
